# Apply updated crypto price/volume figures to columns D (Price) and E (Volume 1h)
# Values in column D are written with a leading apostrophe (quote-prefix) so Excel
# keeps them as literal text (preserving formatting such as trailing zeros and
# thousand-separator dots) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.749.56"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'2.509.49"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'322.36"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'108.46"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.560"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").Value = "'40.33"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'19.51"
$ws.Range("E12").Value = "  +5.32%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "'7.18"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "'2.900.69"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'2.510.39"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'0.851"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'47.664.39"
$ws.Range("D19").Value = "'13.35"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("E22").Value = "  +7.67%  "
$ws.Range("D23").Value = "'71.02"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'247.60"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'25.79"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "'10.20"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "'2.27"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'0.141"
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("D31").Value = "'34.93"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'49.82"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "'0.0785"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'4.71"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'22.57"
$ws.Range("E41").Value = "  +6.01%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'118.22"
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "'2.004.92"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'3.12"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").Value = "'1.81"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'5.17"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "'56.43"
$ws.Range("E51").Value = "  +1.56%  "
